# Guru99BankTestData.xlsx - refresh test-fixture data.
#
# AddNewCustomer (sheet1): the e-mail address typed into the "New Customer"
# form (I2) is refreshed from the old appistoki.com address to a gmail.com
# one, and the selection is left sitting on that cell.
#
# RegisteredCustomer (sheet2): the corresponding "E-mail" column (J2) picks
# up the same new address, and the "Customer ID" (A2) is refreshed to a new
# generated id. Column J is narrower now that the new address is shorter.

$wb = $excel.ActiveWorkbook

$newEmail = "subhash.kiran@gmail.com"
$newCustomerId = "36475"

# --- AddNewCustomer ---------------------------------------------------
$wsAdd = $wb.Worksheets.Item("AddNewCustomer")

$wsAdd.Range("I2").Value = $newEmail

# Leave the cursor on the e-mail cell (was I15).
$wsAdd.Activate()
$wsAdd.Range("I2").Select()

# --- RegisteredCustomer -------------------------------------------------
$wsReg = $wb.Worksheets.Item("RegisteredCustomer")

# A2 ("Customer ID") holds a digit-string that must stay a *text* value
# (it is stored as a shared string, matching every other ID/PIN column on
# this sheet) rather than collapse into a number when assigned straight
# through .Value. Borrow the existing text formatting from a cell that's
# already formatted as Text (AddNewCustomer!G2, the PIN field) via a
# formats-only paste, write through a scratch cell, then paste the result
# into A2 as values-only so A2's own style index ("s=9") never changes.
$wsAdd.Range("G2").Copy()
$wsReg.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$wsReg.Range("Z1").Value = $newCustomerId
$wsReg.Range("Z1").Copy()
$wsReg.Range("A2").PasteSpecial(-4163)   # xlPasteValues
$wsReg.Range("Z1").Clear()

$wsReg.Range("J2").Value = $newEmail

# The shorter e-mail address no longer needs as wide a column.
$wsReg.Columns.Item(10).ColumnWidth = 23.8
